# Regenerate column G ("K") values for rows 2..23 on Sheet1.
# This mirrors the upstream commit "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals" which recalculated the
# strike-count/K column with new data and rewrote the explicit values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 6
    3  = 3
    4  = 5
    5  = 2
    6  = 8
    7  = 5
    8  = 3
    9  = 6
    10 = 2
    11 = 7
    12 = 0
    13 = 5
    14 = 5
    15 = 2
    16 = 3
    17 = 3
    18 = 1
    19 = 3
    20 = 6
    21 = 3
    22 = 4
    23 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
